$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text format before writing so that
# numeric-looking strings (e.g. "325.08", "0.4823") are preserved exactly
# as text rather than being auto-converted into floating point numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.404.87"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.912.40"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  +0.60%  "
$ws.Range("D5").Value = "325.08"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").Value = "0.4823"
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("D8").Value = "0.4072"
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("D9").Value = "0.08223"
$ws.Range("E9").Value = "  +2.53%  "
$ws.Range("D10").Value = "1.021"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("D11").Value = "23.48"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "1.923.63"
$ws.Range("E12").Value = "  +2.67%  "
$ws.Range("D13").Value = "6.053"
$ws.Range("E13").Value = "  +2.20%  "
$ws.Range("D14").Value = "7.208"
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("D15").Value = "90.99"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").Value = "0.06805"
$ws.Range("E16").Value = "  +1.91%  "
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "0.00001037"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").Value = "29.408.91"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "5.633"
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("D23").Value = "11.81"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").Value = "2.174"
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("D25").Value = "2.149.48"
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("D26").Value = "6.625"
$ws.Range("E26").Value = "  +11.63%  "
$ws.Range("D27").Value = "155.95"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("D28").Value = "20.04"
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("D29").Value = "2.104"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D30").Value = "120.38"
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").Value = "0.09556"
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("D33").Value = "5.619"
$ws.Range("E33").Value = "  +5.05%  "
$ws.Range("D34").Value = "3.549"
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").Value = "1.369"
$ws.Range("E35").Value = "  -0.68%  "
$ws.Range("D36").Value = "0.02284"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").Value = "0.06105"
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("D38").Value = "1.176"
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("D39").Value = "0.5971"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("D40").Value = "8.045"
$ws.Range("E40").Value = "  +2.75%  "
$ws.Range("D41").Value = "10.80"
$ws.Range("E41").Value = "  +7.08%  "
$ws.Range("D42").Value = "0.1847"
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("D43").Value = "2.407"
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("D45").Value = "0.07615"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").Value = "12.44"
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("D47").Value = "0.5571"
$ws.Range("E47").Value = "  +1.53%  "
$ws.Range("D48").Value = "1.953"
$ws.Range("E48").Value = "  +1.80%  "
$ws.Range("D49").Value = "117.56"
$ws.Range("E49").Value = "  +3.98%  "
$ws.Range("D50").Value = "2.420"
$ws.Range("E50").Value = "  +3.86%  "
$ws.Range("D51").Value = "72.16"
$ws.Range("E51").Value = "  +0.74%  "

# Restore the original (unstyled / General) cell style now that the text
# values have been committed, so no stray number-format style lingers.
$ws.Range("D2:E51").Style = "Normal"
